$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value2 = "2008_2009"
$ws.Cells.Item(2, 2).Value2 = "inflation"
$ws.Cells.Item(2, 3).Value2 = -0.0236

$ws.Cells.Item(3, 1).Value2 = "2008_2009"
$ws.Cells.Item(3, 2).Value2 = "interest"
$ws.Cells.Item(3, 3).Value2 = 0.07969999999999999

$ws.Cells.Item(4, 1).Value2 = "2008_2009"
$ws.Cells.Item(4, 2).Value2 = "uncertain"
$ws.Cells.Item(4, 3).Value2 = -0.4164

$ws.Cells.Item(5, 1).Value2 = "2008_2009"
$ws.Cells.Item(5, 2).Value2 = "invest"
$ws.Cells.Item(5, 3).Value2 = 0.4366

$ws.Cells.Item(6, 1).Value2 = "2008_2009"
$ws.Cells.Item(6, 2).Value2 = "trade"
$ws.Cells.Item(6, 3).Value2 = 0.8377

$ws.Cells.Item(7, 1).Value2 = "2010_2019"
$ws.Cells.Item(7, 2).Value2 = "uncertain"
$ws.Cells.Item(7, 3).Value2 = 0.06569999999999999

$ws.Cells.Item(8, 1).Value2 = "2010_2019"
$ws.Cells.Item(8, 2).Value2 = "interest"
$ws.Cells.Item(8, 3).Value2 = -0.0998

$ws.Cells.Item(9, 1).Value2 = "2010_2019"
$ws.Cells.Item(9, 2).Value2 = "trade"
$ws.Cells.Item(9, 3).Value2 = -0.09619999999999999

$ws.Cells.Item(10, 1).Value2 = "2010_2019"
$ws.Cells.Item(10, 2).Value2 = "invest"
$ws.Cells.Item(10, 3).Value2 = -0.1881

$ws.Cells.Item(11, 1).Value2 = "2010_2019"
$ws.Cells.Item(11, 2).Value2 = "inflation"
$ws.Cells.Item(11, 3).Value2 = 0.1578

$ws.Cells.Item(12, 1).Value2 = "2020_2021"
$ws.Cells.Item(12, 2).Value2 = "inflation"
$ws.Cells.Item(12, 3).Value2 = -0.2313

$ws.Cells.Item(13, 1).Value2 = "2020_2021"
$ws.Cells.Item(13, 2).Value2 = "interest"
$ws.Cells.Item(13, 3).Value2 = -0.0225

$ws.Cells.Item(14, 1).Value2 = "2020_2021"
$ws.Cells.Item(14, 2).Value2 = "invest"
$ws.Cells.Item(14, 3).Value2 = 0.224

$ws.Cells.Item(15, 1).Value2 = "2020_2021"
$ws.Cells.Item(15, 2).Value2 = "trade"
$ws.Cells.Item(15, 3).Value2 = -0.1777

$ws.Cells.Item(16, 1).Value2 = "2020_2021"
$ws.Cells.Item(16, 2).Value2 = "uncertain"
$ws.Cells.Item(16, 3).Value2 = -0.1577

$ws.Cells.Item(17, 1).Value2 = "2022_2023"
$ws.Cells.Item(17, 2).Value2 = "inflation"
$ws.Cells.Item(17, 3).Value2 = -0.246

$ws.Cells.Item(18, 1).Value2 = "2022_2023"
$ws.Cells.Item(18, 2).Value2 = "interest"
$ws.Cells.Item(18, 3).Value2 = -0.4429

$ws.Cells.Item(19, 1).Value2 = "2022_2023"
$ws.Cells.Item(19, 2).Value2 = "uncertain"
$ws.Cells.Item(19, 3).Value2 = -0.3414

$ws.Cells.Item(20, 1).Value2 = "2022_2023"
$ws.Cells.Item(20, 2).Value2 = "invest"
$ws.Cells.Item(20, 3).Value2 = 0.1738

$ws.Cells.Item(21, 1).Value2 = "2022_2023"
$ws.Cells.Item(21, 2).Value2 = "trade"
$ws.Cells.Item(21, 3).Value2 = -0.7594
